# Generate Report for Handoff
#
# The previous handoff round (uuid 7cb241d1-1a9b-4483-8afd-64024a79a70b, a
# single .md file) is superseded by a new handoff round that includes the
# .md file plus two dependent .png assets. This updates the three report
# sheets (Overview, zh-cn, de-de) to reflect the new round:
#   - the existing "row 2" entries move from the old uuid/hash/timestamps
#     to the new ones
#   - two new rows are appended (one per .png dependency) on every sheet
#   - every "source file" cell keeps its "Ready for handoff" /
#     "Include" / "IsDependency" bookkeeping columns and a hyperlink to
#     the underlying file

$wb = $excel.ActiveWorkbook

$oldUuid = "7cb241d1-1a9b-4483-8afd-64024a79a70b"
$newUuid = "00fe0a52-1577-4dd8-848b-b032c717a9bc"
$oldHash = "bdacf8f57db88224d41e87ef38ac16b62fababf8"
$newHash = "403c03738433c3e18cbe7c3d9d55384d35a71868"

$png1Uuid = "5dd77cd9-f630-4998-a323-31a5b81a677f"
$png2Uuid = "def4bae2-5290-411c-b137-bc5a0187cd7c"
$png1Hash = "ab8083b49c62588c48ce634295cd33e63807a541"
$png2Hash = "d03b41ec984ad4ef73d2e9d8669c9ea957402605"

$mdName = "$newUuid.md"
$png1Name = "$png1Uuid.png"
$png2Name = "$png2Uuid.png"
$png1DepName = "$png1Hash.png"
$png2DepName = "$png2Hash.png"

$overviewDate = "2016-18-17 11:18:13"
$zhXlfName = "$newUuid.$newHash.zh-cn.xlf"
$deXlfName = "$newUuid.$newHash.de-de.xlf"
$zhDate = "2016-03-17 11:18:10"
$deDate = "2016-03-17 11:18:13"
$epoch = "0001-01-01 00:00:00"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/42a048c79dd1abf51d474961a9e416aa5897f135/e2e"
$zhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ff83e6a640e947ff775e453d0252d60b9e63403/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/22dff19d5f94b14648bbe62fbe9b5c5df7d1afa4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

function Set-HyperlinkCell($ws, $addr, $text, $url, $display) {
    $ws.Range($addr).Value = $text
    $ws.Hyperlinks.Add($ws.Range($addr), $url, "", "", $display) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2: roll the existing handoff row from the old uuid to the new one.
Set-HyperlinkCell $ov "A2" $mdName "$mdBase/$mdName" $mdName
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = $overviewDate

# Row 3 / 4: the two new .png dependencies.
Set-HyperlinkCell $ov "A3" $png1Name "$mdBase/$png1Name" $png1Name
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = $overviewDate

Set-HyperlinkCell $ov "A4" $png2Name "$mdBase/$png2Name" $png2Name
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = $overviewDate

# ---------------------------------------------------------------------
# Per-locale sheets ("zh-cn" / "de-de")
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; XlfName = $zhXlfName; Date = $zhDate; Base = $zhBase },
    @{ Sheet = "de-de"; XlfName = $deXlfName; Date = $deDate; Base = $deBase }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Row 2: the primary .md handoff file, rolled to the new uuid/hash.
    Set-HyperlinkCell $ws "A2" $mdName "$mdBase/$mdName" $mdName
    Set-HyperlinkCell $ws "B2" ".md" "$mdBase/$mdName" ".md"
    $ws.Range("C2").Value = "Ready for handoff"
    Set-HyperlinkCell $ws "D2" $loc.XlfName "$($loc.Base)/$($loc.XlfName)" $loc.XlfName
    $ws.Range("E2").Value = $loc.Date
    $ws.Range("H2").Value = $epoch
    $ws.Range("I2").Value = "Include"

    # Row 3: first .png dependency.
    Set-HyperlinkCell $ws "A3" $png1Name "$mdBase/$png1Name" $png1Name
    Set-HyperlinkCell $ws "B3" ".png" "$mdBase/$png1Name" ".png"
    $ws.Range("C3").Value = "Ready for handoff"
    Set-HyperlinkCell $ws "D3" $png1DepName "$($loc.Base)/$png1DepName" $png1DepName
    $ws.Range("E3").Value = $loc.Date
    $ws.Range("H3").Value = $epoch
    $ws.Range("I3").Value = "IsDependency"
    $ws.Range("J3").Value = "e2e\$mdName"

    # Row 4: second .png dependency.
    Set-HyperlinkCell $ws "A4" $png2Name "$mdBase/$png2Name" $png2Name
    Set-HyperlinkCell $ws "B4" ".png" "$mdBase/$png2Name" ".png"
    $ws.Range("C4").Value = "Ready for handoff"
    Set-HyperlinkCell $ws "D4" $png2DepName "$($loc.Base)/$png2DepName" $png2DepName
    $ws.Range("E4").Value = $loc.Date
    $ws.Range("H4").Value = $epoch
    $ws.Range("I4").Value = "IsDependency"
    $ws.Range("J4").Value = "e2e\$mdName"
}

Write-Host "Report regenerated for handoff round $newUuid"
